$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 570. This shifts the existing row 570
# ("Calameño"/"Cuarta", 2021-01-15) and everything below it down by one
# (old row 681 becomes row 682), matching the updated dimension A1:R682.
$ws.Rows.Item(570).Insert()

# Populate the newly inserted row 570 with the new "Tuna"/"Segunda" record
# (columns A, B, C, E, F, G, R keep the same constant values shared by all
# rows in this block).
$ws.Range("A570").Value = 6
$ws.Range("B570").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C570").Value = "Metropolitana"
$ws.Range("D570").Value = 44508
$ws.Range("E570").Value = 13
$ws.Range("F570").Value = 100112027
$ws.Range("G570").Value = "Melón"
$ws.Range("H570").Value = "Tuna"
$ws.Range("I570").Value = "Segunda"
$ws.Range("J570").Value = 340
$ws.Range("K570").Value = 27000
$ws.Range("L570").Value = 28000
$ws.Range("M570").Value = 27559
$ws.Range("N570").Value = "$/caja 24 unidades"
$ws.Range("O570").Value = "Provincia de Copiapó"
$ws.Range("P570").Value = 1148
$ws.Range("Q570").Value = 24
$ws.Range("R570").Value = "Hortaliza"
